# Update the "Users to Delete" sheet with new Database ID (GUID) values
$wb = $excel.ActiveWorkbook

$wsUsersToDelete = $wb.Worksheets.Item("Users to Delete")
$wsUsersToDelete.Range("D2").Value = "8f2d0c08-f7ad-4935-a6e5-c52dde1be5e0"
$wsUsersToDelete.Range("D3").Value = "1d2c95f6-a889-49b9-9bcc-3908c7ca40c1"
$wsUsersToDelete.Range("D4").Value = "7b0f8a3b-b7a0-4b9f-a880-a945f60ee21e"
$wsUsersToDelete.Range("D5").Value = "fd395b8e-1bbb-4a51-829b-27a6beb7f4f3"

# Update the "Summary" sheet's Report Generated timestamp
$wsSummary = $wb.Worksheets.Item("Summary")
$wsSummary.Range("B6").Value = "10/28/2025, 8:17:23 PM"
